# Revert the "fild" merge: restore the original team/role assignments for the
# intern roster (Sheet1) and the pre-merge Normal/Hyperlink font (Tahoma).
#
# This reverts a merge that had rewritten six interns' department / team /
# mentor / role columns (C:G) to a different project line-up (CRM / Mobile /
# Marketing / BMs ...). Columns A,B,H,I,J,K,N,O,P,Q,R (id, access request,
# prefix, name, nickname, start date, phone, email, institute, faculty,
# major) are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Wasaya Phaisan / Fogus) ---------------------------------------
$ws.Range("C2").Value = "Accounting"
$ws.Range("D2").Value = $null
$ws.Range("E2").Value = "Acc Team A"
$ws.Range("F2").Value = "กิตติธร ปรีดาอัครกุล"
$ws.Range("G2").Value = "Accounting"

# --- Row 3 (Paweena Luekham / Paula) --------------------------------------
$ws.Range("C3").Value = "Accounting"
$ws.Range("D3").Value = $null
$ws.Range("E3").Value = "Acc Team A"
$ws.Range("F3").Value = "กิตติธร ปรีดาอัครกุล"
$ws.Range("G3").Value = "Accounting"

# --- Row 4 (Preeyawadee Kiatbowornsakul / Ren) ----------------------------
$ws.Range("C4").Value = "BMS"
$ws.Range("D4").Value = "Developer"
$ws.Range("E4").Value = "BMS Team B"
$ws.Range("F4").Value = "พัฒน์ สุพรรณภาคิน"
$ws.Range("G4").Value = "Business Analyst"

# --- Row 5 (Pitipoom Watthanasakmontri / Peem) ----------------------------
$ws.Range("C5").Value = "Government Project"
$ws.Range("D5").Value = $null
$ws.Range("E5").Value = "Gov Team A"
$ws.Range("F5").Value = "ปริญญา ศิลาดี"
$ws.Range("G5").Value = "Business Analyst"

# --- Row 6 (Jak Rattanaprasert / Zee) -------------------------------------
$ws.Range("C6").Value = "Interactive Media"
$ws.Range("D6").Value = $null
$ws.Range("E6").Value = "Mobile A"
$ws.Range("F6").Value = "ปฐมพร ภูพาณิชย์"
$ws.Range("G6").Value = "Developer"

# --- Row 7 (Phanida Thamwapee / Ice) --------------------------------------
$ws.Range("C7").Value = "BMS"
$ws.Range("D7").Value = "Developer"
$ws.Range("E7").Value = "BMS Team A"
$ws.Range("F7").Value = "พัฒน์ สุพรรณภาคิน"
$ws.Range("G7").Value = "Developer"

# --- Restore the pre-merge body/hyperlink font (Calibri -> Tahoma) --------
$used = $ws.UsedRange
$used.Font.Name = "Tahoma"
$ws.Hyperlinks.Item(1).Range.Font.Name = "Tahoma"
$ws.Hyperlinks.Item(2).Range.Font.Name = "Tahoma"
$ws.Hyperlinks.Item(3).Range.Font.Name = "Tahoma"
$ws.Hyperlinks.Item(4).Range.Font.Name = "Tahoma"
$ws.Hyperlinks.Item(5).Range.Font.Name = "Tahoma"
$ws.Hyperlinks.Item(6).Range.Font.Name = "Tahoma"

# --- Clear the stray selection / scroll position left on the sheet -------
$ws.Range("A1").Select()
